$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Education" / "Nothing; don't know; empty" labels (rows 23 & 24)
$ws.Cells.Item(23, 1).Value = "Education"
$ws.Cells.Item(24, 1).Value = "Nothing; don't know; empty"

# Re-run RU (Russia, column L) 1001 without crop -> refreshed percentages
$ws.Cells.Item(2, 2).Value = 0.288602254068036
$ws.Cells.Item(2, 12).Value = 0.256470409182717
$ws.Cells.Item(3, 2).Value = 0.121669824607341
$ws.Cells.Item(3, 12).Value = 0.0837834685235709
$ws.Cells.Item(4, 2).Value = 0.120833344367493
$ws.Cells.Item(4, 12).Value = 0.0742699457960542
$ws.Cells.Item(5, 2).Value = 0.120389085503789
$ws.Cells.Item(5, 12).Value = 0.267963612801303
$ws.Cells.Item(6, 2).Value = 0.120184153459867
$ws.Cells.Item(6, 12).Value = 0.285935356678167
$ws.Cells.Item(7, 2).Value = 0.0952935613072265
$ws.Cells.Item(7, 12).Value = 0.124618910276095
$ws.Cells.Item(8, 2).Value = 0.0794735282239128
$ws.Cells.Item(8, 12).Value = 0.0083177492120678
$ws.Cells.Item(9, 2).Value = 0.0735093679559291
$ws.Cells.Item(10, 2).Value = 0.0698768763127664
$ws.Cells.Item(10, 12).Value = 0.111081482769912
$ws.Cells.Item(11, 2).Value = 0.0635424288770414
$ws.Cells.Item(11, 12).Value = 0.0412759090274652
$ws.Cells.Item(12, 2).Value = 0.0604260895086438
$ws.Cells.Item(12, 12).Value = 0.0139423840542105
$ws.Cells.Item(13, 2).Value = 0.0539494646923607
$ws.Cells.Item(14, 2).Value = 0.0539186966139836
$ws.Cells.Item(14, 12).Value = 0.0818131194776077
$ws.Cells.Item(15, 2).Value = 0.0508992953124374
$ws.Cells.Item(16, 2).Value = 0.042695078995228
$ws.Cells.Item(16, 12).Value = 0.00948397951163939
$ws.Cells.Item(17, 2).Value = 0.0389291469046823
$ws.Cells.Item(17, 12).Value = 0.0571120374959608
$ws.Cells.Item(18, 2).Value = 0.0302860920786351
$ws.Cells.Item(19, 2).Value = 0.0297938356166113
$ws.Cells.Item(19, 12).Value = 0.0279987225267749
$ws.Cells.Item(20, 2).Value = 0.0292523790861179
$ws.Cells.Item(21, 2).Value = 0.0291765410820688
$ws.Cells.Item(22, 2).Value = 0.0263602206785098
$ws.Cells.Item(22, 12).Value = 0.0145849886052593
$ws.Cells.Item(23, 2).Value = 0.020720803874429
$ws.Cells.Item(23, 3).Value = 0.0196837537441796
$ws.Cells.Item(23, 4).Value = 0.0182524677661447
$ws.Cells.Item(23, 5).Value = 0.0190783171773221
$ws.Cells.Item(23, 6).Value = 0.0088034747635748
$ws.Cells.Item(23, 8).Value = 0.0610540047534814
$ws.Cells.Item(23, 9).Value = 0.0227761038083329
$ws.Cells.Item(23, 10).Value = 0.00615537143247198
$ws.Cells.Item(23, 11).Value = 0.0160384531819767
$ws.Cells.Item(23, 12).Value = 0.0287909007211753
$ws.Cells.Item(23, 13).Value = 0.0613354638950244
$ws.Cells.Item(23, 14).Value = 0.0176125532277882
$ws.Cells.Item(24, 2).Value = 0.0202796245140543
$ws.Cells.Item(24, 3).Value = 0.0224339874020708
$ws.Cells.Item(24, 4).Value = 0.032429789172551
$ws.Cells.Item(24, 5).Value = 0.0354546001411635
$ws.Cells.Item(24, 6).Value = 0.0217770349278578
$ws.Cells.Item(24, 8).Value = 0.0395234109976004
$ws.Cells.Item(24, 9).Value = 0.00390389155818829
$ws.Cells.Item(24, 10).Value = 0.0151668384924855
$ws.Cells.Item(24, 11).Value = 0.0640268258313086
$ws.Cells.Item(24, 12).Value = 0.00427350427350427
$ws.Cells.Item(24, 13).Value = 0.00461237333330672
$ws.Cells.Item(24, 14).Value = 0.00485375235334099
$ws.Cells.Item(25, 2).Value = 0.0173330250156156
$ws.Cells.Item(25, 12).Value = 0.068041723121496
$ws.Cells.Item(26, 2).Value = 0.0170562849683906
$ws.Cells.Item(27, 2).Value = 0.0168462744348631
$ws.Cells.Item(27, 12).Value = 0.0571281230595419
$ws.Cells.Item(28, 2).Value = 0.0153164453736903
$ws.Cells.Item(28, 12).Value = 0.00948397951163939
$ws.Cells.Item(29, 2).Value = 0.0119642459635955
$ws.Cells.Item(29, 3).Value = 0.00839210229336858
$ws.Cells.Item(29, 12).Value = 0.00805580096232858
$ws.Cells.Item(30, 2).Value = 0.0111426227379494
$ws.Cells.Item(31, 2).Value = 0.00973820469847729
$ws.Cells.Item(32, 2).Value = 0.00925328230690155
$ws.Cells.Item(32, 12).Value = 0.0441176856446852
$ws.Cells.Item(33, 2).Value = 0.0087253775238692
$ws.Cells.Item(34, 2).Value = 0.00847087921126582
$ws.Cells.Item(34, 12).Value = 0.00213675213675214
$ws.Cells.Item(35, 2).Value = 0.00622786349050728
$ws.Cells.Item(36, 2).Value = 0.00535683716861164
$ws.Cells.Item(37, 2).Value = 0.00509441459479364
$ws.Cells.Item(37, 12).Value = 0.0233795998838649
$ws.Cells.Item(38, 2).Value = 0.00509169877264681
$ws.Cells.Item(38, 12).Value = 0.00805580096232858
$ws.Cells.Item(39, 2).Value = 0.00413381453197999
$ws.Cells.Item(40, 2).Value = 0.00201383583398161
$ws.Cells.Item(41, 2).Value = 0.00187011792770167
$ws.Cells.Item(42, 2).Value = 0.00183177056893757
$ws.Cells.Item(43, 2).Value = 0.00133715087982487
$ws.Cells.Item(44, 2).Value = 0.00123053212590171
$ws.Cells.Item(45, 2).Value = 0.0010612055022385
$ws.Cells.Item(46, 2).Value = 0.0010612055022385
$ws.Cells.Item(47, 2).Value = 0.000883290542613318
$ws.Cells.Item(48, 2).Value = 0.000792106314689323
